$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 7937
$ws.Range("C3").Value = 13600
$ws.Range("D3").Value = 4035
$ws.Range("E3").Value = 8249
$ws.Range("F3").Value = 12500
$ws.Range("G3").Value = 12900

$ws.Range("B4").Value = 1040.187392
$ws.Range("C4").Value = 1777.33632
$ws.Range("D4").Value = 528.482304
$ws.Range("E4").Value = 1081.081856
$ws.Range("F4").Value = 1636.827136
$ws.Range("G4").Value = 1690.304512

$ws.Range("B5").Value = 125.24
$ws.Range("C5").Value = 144.85
$ws.Range("D5").Value = 911.16
$ws.Range("E5").Value = 805.22
$ws.Range("F5").Value = 1191.82
$ws.Range("G5").Value = 2373.09

$ws.Range("B6").Value = 217
$ws.Range("C6").Value = 233
$ws.Range("D6").Value = 1827
$ws.Range("E6").Value = 1893
$ws.Range("F6").Value = 2704
$ws.Range("G6").Value = 5735

$ws.Range("B7").Value = 322
$ws.Range("C7").Value = 338
$ws.Range("D7").Value = 2057
$ws.Range("E7").Value = 2114
$ws.Range("F7").Value = 3752
$ws.Range("G7").Value = 7308

$ws.Range("B12").Value = 5879
$ws.Range("C12").Value = 46100
$ws.Range("D12").Value = 99100
$ws.Range("E12").Value = 143000
$ws.Range("F12").Value = 63300
$ws.Range("G12").Value = 436000

$ws.Range("B13").Value = 23.068672
$ws.Range("C13").Value = 188.74368
$ws.Range("D13").Value = 405.798912
$ws.Range("E13").Value = 587.2025599999999
$ws.Range("F13").Value = 258.998272
$ws.Range("G13").Value = 1784.676352

$ws.Range("B14").Value = 168.87664
$ws.Range("C14").Value = 42.28108
$ws.Range("D14").Value = 38.28733
$ws.Range("E14").Value = 41.92276
$ws.Range("F14").Value = 243.65374
$ws.Range("G14").Value = 64.05472

$ws.Range("B15").Value = 1400.832
$ws.Range("C15").Value = 71.16800000000001
$ws.Range("D15").Value = 71.16800000000001
$ws.Range("E15").Value = 90.624
$ws.Range("F15").Value = 1499.136
$ws.Range("G15").Value = 179.2

$ws.Range("B16").Value = 1548.288
$ws.Range("C16").Value = 77.312
$ws.Range("D16").Value = 82.432
$ws.Range("E16").Value = 110.08
$ws.Range("F16").Value = 1744.896
$ws.Range("G16").Value = 259.072

$ws.Range("B21").Value = 9225
$ws.Range("C21").Value = 16000
$ws.Range("D21").Value = 16800
$ws.Range("E21").Value = 19100
$ws.Range("F21").Value = 20700
$ws.Range("G21").Value = 20500

$ws.Range("B22").Value = 1209.008128
$ws.Range("C22").Value = 2097.152
$ws.Range("D22").Value = 2199.912448
$ws.Range("E22").Value = 2502.950912
$ws.Range("F22").Value = 2708.471808
$ws.Range("G22").Value = 2681.208832

$ws.Range("B23").Value = 52.88
$ws.Range("C23").Value = 33.66
$ws.Range("D23").Value = 38.03
$ws.Range("E23").Value = 64
$ws.Range("F23").Value = 89.31
$ws.Range("G23").Value = 144.88

$ws.Range("B24").Value = 114
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 70
$ws.Range("E24").Value = 149
$ws.Range("F24").Value = 133
$ws.Range("G24").Value = 200

$ws.Range("B25").Value = 922
$ws.Range("C25").Value = 186
$ws.Range("D25").Value = 196
$ws.Range("E25").Value = 243
$ws.Range("F25").Value = 388
$ws.Range("G25").Value = 971

$ws.Range("B30").Value = 328000
$ws.Range("C30").Value = 415000
$ws.Range("D30").Value = 537000
$ws.Range("E30").Value = 586000
$ws.Range("F30").Value = 558000
$ws.Range("G30").Value = 675000

$ws.Range("B31").Value = 1342.17728
$ws.Range("C31").Value = 1698.69312
$ws.Range("D31").Value = 2199.912448
$ws.Range("E31").Value = 2402.287616
$ws.Range("F31").Value = 2284.847104
$ws.Range("G31").Value = 2765.094912

$ws.Range("B32").Value = 1.10329
$ws.Range("C32").Value = 1.19733
$ws.Range("D32").Value = 1.3056
$ws.Range("E32").Value = 1.68041
$ws.Range("F32").Value = 2.76503
$ws.Range("G32").Value = 4.58971

$ws.Range("B33").Value = 1.688
$ws.Range("C33").Value = 1.88
$ws.Range("D33").Value = 1.976
$ws.Range("E33").Value = 2.928
$ws.Range("F33").Value = 4.832
$ws.Range("G33").Value = 6.496

$ws.Range("B34").Value = 6.368
$ws.Range("C34").Value = 6.688
$ws.Range("D34").Value = 6.432
$ws.Range("E34").Value = 7.072
$ws.Range("F34").Value = 9.92
$ws.Range("G34").Value = 10.944

$ws.Range("B39").Value = 17700
$ws.Range("C39").Value = 32000
$ws.Range("D39").Value = 25900
$ws.Range("E39").Value = 8677
$ws.Range("F39").Value = 8155
$ws.Range("G39").Value = 10400

$ws.Range("B40").Value = 2314.207232
$ws.Range("C40").Value = 4194.304
$ws.Range("D40").Value = 3398.434816
$ws.Range("E40").Value = 1137.70496
$ws.Range("F40").Value = 1068.498944
$ws.Range("G40").Value = 1364.197376

$ws.Range("B41").Value = 55.6
$ws.Range("C41").Value = 60.78
$ws.Range("D41").Value = 153.15
$ws.Range("E41").Value = 902.46
$ws.Range("F41").Value = 1950.27
$ws.Range("G41").Value = 3063.62

$ws.Range("B42").Value = 151
$ws.Range("C42").Value = 161
$ws.Range("D42").Value = 383
$ws.Range("E42").Value = 3261
$ws.Range("F42").Value = 7701
$ws.Range("G42").Value = 13829

$ws.Range("B43").Value = 165
$ws.Range("C43").Value = 206
$ws.Range("D43").Value = 469
$ws.Range("E43").Value = 4555
$ws.Range("F43").Value = 9634
$ws.Range("G43").Value = 16188

$ws.Range("B48").Value = 155000
$ws.Range("C48").Value = 169000
$ws.Range("D48").Value = 263000
$ws.Range("E48").Value = 277000
$ws.Range("F48").Value = 249000
$ws.Range("G48").Value = 562000

$ws.Range("B49").Value = 636.485632
$ws.Range("C49").Value = 692.06016
$ws.Range("D49").Value = 1077.936128
$ws.Range("E49").Value = 1133.510656
$ws.Range("F49").Value = 1019.215872
$ws.Range("G49").Value = 2301.62432

$ws.Range("B50").Value = 6.02574
$ws.Range("C50").Value = 9.709899999999999
$ws.Range("D50").Value = 14.51409
$ws.Range("E50").Value = 27.85442
$ws.Range("F50").Value = 63.31621
$ws.Range("G50").Value = 56.08361

$ws.Range("B51").Value = 1.688
$ws.Range("C51").Value = 1.832
$ws.Range("D51").Value = 2.008
$ws.Range("E51").Value = 2.288
$ws.Range("F51").Value = 2.8
$ws.Range("G51").Value = 2.768

$ws.Range("B52").Value = 84.48
$ws.Range("C52").Value = 113.152
$ws.Range("D52").Value = 444.416
$ws.Range("E52").Value = 1122.304
$ws.Range("F52").Value = 2932.736
$ws.Range("G52").Value = 1941.504

$ws.Range("B57").Value = 8677
$ws.Range("C57").Value = 11500
$ws.Range("D57").Value = 12300
$ws.Range("E57").Value = 13700
$ws.Range("F57").Value = 13600
$ws.Range("G57").Value = 12600

$ws.Range("B58").Value = 1137.70496
$ws.Range("C58").Value = 1507.852288
$ws.Range("D58").Value = 1612.709888
$ws.Range("E58").Value = 1795.162112
$ws.Range("F58").Value = 1784.676352
$ws.Range("G58").Value = 1654.652928

$ws.Range("B59").Value = 27.19
$ws.Range("C59").Value = 28.1
$ws.Range("D59").Value = 33.07
$ws.Range("E59").Value = 42.57
$ws.Range("F59").Value = 83.59999999999999
$ws.Range("G59").Value = 142.68

$ws.Range("B60").Value = 30
$ws.Range("C60").Value = 31
$ws.Range("D60").Value = 37
$ws.Range("E60").Value = 60
$ws.Range("F60").Value = 131
$ws.Range("G60").Value = 186

$ws.Range("B61").Value = 110
$ws.Range("C61").Value = 115
$ws.Range("D61").Value = 116
$ws.Range("E61").Value = 126
$ws.Range("F61").Value = 184
$ws.Range("G61").Value = 223

$ws.Range("B66").Value = 241000
$ws.Range("C66").Value = 338000
$ws.Range("D66").Value = 405000
$ws.Range("E66").Value = 438000
$ws.Range("F66").Value = 423000
$ws.Range("G66").Value = 402000

$ws.Range("B67").Value = 986.710016
$ws.Range("C67").Value = 1384.12032
$ws.Range("D67").Value = 1656.75008
$ws.Range("E67").Value = 1793.06496
$ws.Range("F67").Value = 1733.296128
$ws.Range("G67").Value = 1647.312896

$ws.Range("B68").Value = 1.01687
$ws.Range("C68").Value = 1.04411
$ws.Range("D68").Value = 1.12142
$ws.Range("E68").Value = 1.3743
$ws.Range("F68").Value = 2.51215
$ws.Range("G68").Value = 4.24955

$ws.Range("B69").Value = 3.088
$ws.Range("C69").Value = 1.8
$ws.Range("D69").Value = 1.896
$ws.Range("E69").Value = 2.736
$ws.Range("F69").Value = 4.448
$ws.Range("G69").Value = 6.048

$ws.Range("B70").Value = 4.256
$ws.Range("C70").Value = 4.32
$ws.Range("D70").Value = 4.448
$ws.Range("E70").Value = 5.344
$ws.Range("F70").Value = 9.152000000000001
$ws.Range("G70").Value = 13.376
